$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting: bold navy-blue font for the H:J "real time / done / postponed" columns ---
$navyBlue = 6299648   # RGB(0,32,96) = 0 + 32*256 + 96*65536 => navy blue FF002060

$rngH = $ws.Range("H4:H21")
$rngH.Font.Bold = $true
$rngH.Font.Color = $navyBlue

$rngI = $ws.Range("I4:I21")
$rngI.Font.Bold = $true
$rngI.Font.Color = $navyBlue
$rngI.HorizontalAlignment = -4108   # xlCenter
$rngI.VerticalAlignment = -4108     # xlCenter

$rngJ = $ws.Range("J4:J21")
$rngJ.Font.Bold = $true
$rngJ.Font.Color = $navyBlue
$rngJ.VerticalAlignment = -4108     # xlCenter

# --- Data entry: Prefab Manager task (row 5) logged as 60 minutes real time and marked done ---
$ws.Range("H5").Value = 60
$ws.Range("I5").Value = "Oui"

# --- Stray formatting leftover on M11 (strikethrough, no value) ---
$ws.Range("M11").Font.Strikethrough = $true

# --- Move the active selection to M16 ---
[void]$ws.Range("M16").Select()
